$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $value) {
    $range = $ws.Range($addr)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows 2-48
Set-TextCell $ws "D2" "42.980.22"
Set-TextCell $ws "E2" "  -0.26%  "
Set-TextCell $ws "D3" "2.307.55"
Set-TextCell $ws "E3" "  +0.13%  "
Set-TextCell $ws "D4" "1.00"
Set-TextCell $ws "E4" "  +0.04%  "
Set-TextCell $ws "D5" "306.49"
Set-TextCell $ws "E5" "  +2.11%  "
Set-TextCell $ws "D6" "96.85"
Set-TextCell $ws "E6" "  -1.26%  "
Set-TextCell $ws "D7" "0.510"
Set-TextCell $ws "E7" "  -1.79%  "
Set-TextCell $ws "D8" "1.00"
Set-TextCell $ws "E8" "  -0.02%  "
Set-TextCell $ws "D9" "0.504"
Set-TextCell $ws "E9" "  -2.40%  "
Set-TextCell $ws "D10" "35.54"
Set-TextCell $ws "E10" "  -1.82%  "
Set-TextCell $ws "E11" "  +0.13%  "
Set-TextCell $ws "E12" "  +3.95%  "
Set-TextCell $ws "E13" "  +1.22%  "
Set-TextCell $ws "D14" "6.78"
Set-TextCell $ws "E14" "  -1.54%  "
Set-TextCell $ws "D15" "2.657.91"
Set-TextCell $ws "E15" "  -0.11%  "
Set-TextCell $ws "D16" "2.296.46"
Set-TextCell $ws "E16" "  -0.21%  "
Set-TextCell $ws "D17" "0.784"
Set-TextCell $ws "E17" "  -0.58%  "
Set-TextCell $ws "D18" "42.891.16"
Set-TextCell $ws "E18" "  -0.20%  "
Set-TextCell $ws "D19" "12.99"
Set-TextCell $ws "E19" "  +2.01%  "
Set-TextCell $ws "D20" "0.0₃0899"
Set-TextCell $ws "E20" "  -1.35%  "
Set-TextCell $ws "D21" "6.05"
Set-TextCell $ws "E21" "  -1.69%  "
Set-TextCell $ws "D22" "67.44"
Set-TextCell $ws "E22" "  -1.79%  "
Set-TextCell $ws "D23" "236.64"
Set-TextCell $ws "E23" "  -0.59%  "
Set-TextCell $ws "D24" "2.14"
Set-TextCell $ws "E24" "  -1.46%  "
Set-TextCell $ws "E25" "  +1.06%  "
Set-TextCell $ws "E26" "  +0.14%  "
Set-TextCell $ws "E27" "  +0.03%  "
Set-TextCell $ws "E28" "  +1.22%  "
Set-TextCell $ws "E29" "  +6.44%  "
Set-TextCell $ws "D30" "166.80"
Set-TextCell $ws "E30" "  +1.47%  "
Set-TextCell $ws "D31" "9.09"
Set-TextCell $ws "E31" "  -0.64%  "
Set-TextCell $ws "D32" "33.23"
Set-TextCell $ws "E32" "  +0.42%  "
Set-TextCell $ws "E33" "  +0.10%  "
Set-TextCell $ws "D34" "4.79"
Set-TextCell $ws "E34" "  -0.70%  "
Set-TextCell $ws "D35" "5.00"
Set-TextCell $ws "E35" "  -2.29%  "
Set-TextCell $ws "D36" "17.81"
Set-TextCell $ws "E36" "  -1.68%  "
Set-TextCell $ws "E37" "  -0.99%  "
Set-TextCell $ws "E38" "  -0.56%  "
Set-TextCell $ws "E39" "  -0.84%  "
Set-TextCell $ws "E40" "  -1.24%  "
Set-TextCell $ws "E41" "  -0.91%  "
Set-TextCell $ws "E42" "  -2.79%  "
Set-TextCell $ws "D43" "2.012.47"
Set-TextCell $ws "E43" "  -0.36%  "
Set-TextCell $ws "D44" "0.0281"
Set-TextCell $ws "E44" "  -2.08%  "
Set-TextCell $ws "D45" "18.20"
Set-TextCell $ws "E45" "  +3.99%  "
Set-TextCell $ws "D46" "10.05"
Set-TextCell $ws "E46" "  -3.33%  "
Set-TextCell $ws "D47" "2.06"
Set-TextCell $ws "E47" "  -6.92%  "
Set-TextCell $ws "D48" "2.80"
Set-TextCell $ws "E48" "  -1.58%  "

# Rows 49-51 shifted: new coin (HuobiToken) inserted, BitcoinSV dropped
Set-TextCell $ws "B49" "HuobiToken"
Set-TextCell $ws "C49" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D49" "2.89"
Set-TextCell $ws "E49" "  +11.07%  "

Set-TextCell $ws "B50" "MultiversX"
Set-TextCell $ws "C50" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws "D50" "53.97"
Set-TextCell $ws "E50" "  -0.47%  "

Set-TextCell $ws "B51" "RocketPoolETH"
Set-TextCell $ws "C51" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws "D51" "2.525.90"
Set-TextCell $ws "E51" "  -0.18%  "
